$d = $word.ActiveDocument

# 1. M2DocEvaluator.caseLet line: 1228 -> 1251
$d.Content.Find.Execute("caseLet(M2DocEvaluator.java:1228)", $true, $false, $false, $false, $false, $true, 1, $false, "caseLet(M2DocEvaluator.java:1251)", 2)

# 2. M2DocEvaluator.doSwitch line: 1216 -> 1239 (3 occurrences, replace all)
$d.Content.Find.Execute("doSwitch(M2DocEvaluator.java:1216)", $true, $false, $false, $false, $false, $true, 1, $false, "doSwitch(M2DocEvaluator.java:1239)", 2)

# 3. M2DocEvaluator.caseBlock line: 1425 -> 1464
$d.Content.Find.Execute("caseBlock(M2DocEvaluator.java:1425)", $true, $false, $false, $false, $false, $true, 1, $false, "caseBlock(M2DocEvaluator.java:1464)", 2)

# 4. M2DocEvaluator.caseDocumentTemplate line: 287 -> 296
$d.Content.Find.Execute("caseDocumentTemplate(M2DocEvaluator.java:287)", $true, $false, $false, $false, $false, $true, 1, $false, "caseDocumentTemplate(M2DocEvaluator.java:296)", 2)

# 5. M2DocEvaluator.generate line: 276 -> 281
$d.Content.Find.Execute("generate(M2DocEvaluator.java:276)", $true, $false, $false, $false, $false, $true, 1, $false, "generate(M2DocEvaluator.java:281)", 2)

# 6. M2DocUtils.generate line: 694 -> 805
$d.Content.Find.Execute("M2DocUtils.generate(M2DocUtils.java:694)", $true, $false, $false, $false, $false, $true, 1, $false, "M2DocUtils.generate(M2DocUtils.java:805)", 2)

# 7. AbstractTemplatesTestSuite.prepareoutputAndGenerate line: 480 -> 511
$d.Content.Find.Execute("prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:480)", $true, $false, $false, $false, $false, $true, 1, $false, "prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:511)", 2)

# 8. AbstractTemplatesTestSuite.generation line: 389 -> 420
$d.Content.Find.Execute("generation(AbstractTemplatesTestSuite.java:389)", $true, $false, $false, $false, $false, $true, 1, $false, "generation(AbstractTemplatesTestSuite.java:420)", 2)

# 9. Insert new stack trace line before second occurrence of RunAfters.evaluate(RunAfters.java:27)
$oldChunk = [char]9 + "at org.junit.runners.ParentRunner" + [char]36 + "2.evaluate(ParentRunner.java:268)" + [char]10 + [char]9 + "at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)" + [char]10 + [char]9 + "at org.junit.runners.ParentRunner.run(ParentRunner.java:363)" + [char]10 + [char]9 + "at org.junit.runners.Suite.runChild(Suite.java:128)" + [char]10 + [char]9 + "at org.junit.runners.Suite.runChild(Suite.java:27)" + [char]10 + [char]9 + "at org.junit.runners.ParentRunner" + [char]36 + "3.run(ParentRunner.java:290)"

$newChunk = [char]9 + "at org.junit.runners.ParentRunner" + [char]36 + "2.evaluate(ParentRunner.java:268)" + [char]10 + [char]9 + "at org.junit.internal.runners.statements.RunBefores.evaluate(RunBefores.java:26)" + [char]10 + [char]9 + "at org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)" + [char]10 + [char]9 + "at org.junit.runners.ParentRunner.run(ParentRunner.java:363)" + [char]10 + [char]9 + "at org.junit.runners.Suite.runChild(Suite.java:128)" + [char]10 + [char]9 + "at org.junit.runners.Suite.runChild(Suite.java:27)" + [char]10 + [char]9 + "at org.junit.runners.ParentRunner" + [char]36 + "3.run(ParentRunner.java:290)"

$d.Content.Find.Execute($oldChunk, $true, $false, $false, $false, $false, $true, 1, $false, $newChunk, 2)

Write-Output "Done"
